# Apply "Holden scheme" update to the simulation worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the duplicate trailing block of columns U:AD (they mirrored K:T)
#    so the sheet shrinks back down to columns A:T.
$ws.Range("U1:AD1").EntireColumn.Delete()

# 2) Row 2 holds the column headers for C:T. Re-order the "[h, k, l]"
#    group and re-point the rest at the same label set (values only move,
#    the trailing duplicate set from step 1 already dropped off).
$headerVals = @("[3, 2, 1]", "[1, 1, 0]", "[2, 2, 2]", "[3, 1, 0]", "[2, 2, 0]", "[2, 0, 0]", "[2, 1, 1]", "[4, 0, 0]", "1Pair-A", "1Pair-B", "2Pairs-A", "2Pairs-B", "3Pairs-A", "3Pairs-B", "3Pairs-C", "4Pairs", "5A4F", "MaxUnique")
for ($i = 0; $i -lt $headerVals.Length; $i++) {
    $col = 3 + $i
    $ws.Cells.Item(2, $col).Value = $headerVals[$i]
}

# 3) Rows 16-19 were the "HexGrid-90degTilt*" entries; rename them to the
#    new "Holden*" scheme labels (the underlying numbers in C:T stay 1).
$ws.Cells.Item(16, 2).Value = "Holden2.5"
$ws.Cells.Item(17, 2).Value = "Holden5"
$ws.Cells.Item(18, 2).Value = "Holden10"
$ws.Cells.Item(19, 2).Value = "Holden15"

# 4) Append four new rows (20-23) re-using the "HexGrid-90degTilt*" labels
#    that moved off of rows 16-19, each filled with 1s across C:T like the
#    existing data rows.
$newRows = @(
    @{ Row = 20; A = 18; Label = "HexGrid-90degTilt2.5degRes" },
    @{ Row = 21; A = 19; Label = "HexGrid-90degTilt5degRes" },
    @{ Row = 22; A = 20; Label = "HexGrid-90degTilt10degRes" },
    @{ Row = 23; A = 21; Label = "HexGrid-90degTilt15degRes" }
)

# Column A uses the bordered/bold "header" style (style index 1 in the
# original file) - copy that formatting from an existing A-column cell
# down onto the new rows instead of trying to set it from scratch.
$ws.Range("A19").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = $nr.A
    $ws.Cells.Item($r, 2).Value = $nr.Label
    for ($col = 3; $col -le 20; $col++) {
        $ws.Cells.Item($r, $col).Value = 1
    }
}
